$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

    $ws.Range("B2").Value = 28.49229048257739
    $ws.Range("C2").Value = 10.57254853506121
    $ws.Range("D2").Value = 3.830039077387625
    $ws.Range("F2").Value = 56.24765701318223
    $ws.Range("G2").Value = 3.791572029002654
    $ws.Range("J2").Value = 9.821912057493767
    $ws.Range("L2").Value = 11.70496252890569
    $ws.Range("M2").Value = 21.9996927765327
    $ws.Range("N2").Value = 23.53048010724265
    $ws.Range("B3").Value = 28.24294725512831
    $ws.Range("C3").Value = 10.19742042821624
    $ws.Range("D3").Value = 3.76759115399786
    $ws.Range("F3").Value = 56.153901398131
    $ws.Range("G3").Value = 3.795981986361284
    $ws.Range("J3").Value = 9.829207313864083
    $ws.Range("L3").Value = 11.72431416743099
    $ws.Range("M3").Value = 21.97347256535217
    $ws.Range("N3").Value = 23.58581786079055
    $ws.Range("B4").Value = 28.09680786055588
    $ws.Range("C4").Value = 9.963695142295386
    $ws.Range("D4").Value = 3.728138905517848
    $ws.Range("F4").Value = 56.10886182582507
    $ws.Range("G4").Value = 3.798829396249521
    $ws.Range("J4").Value = 9.83392495012912
    $ws.Range("L4").Value = 11.7376838585659
    $ws.Range("M4").Value = 21.96193486391492
    $ws.Range("N4").Value = 23.6217764214887
    $ws.Range("B5").Value = 28.03906233666212
    $ws.Range("C5").Value = 9.86776813456725
    $ws.Range("D5").Value = 3.711789051060111
    $ws.Range("F5").Value = 56.09366072387819
    $ws.Range("G5").Value = 3.800025002313178
    $ws.Range("J5").Value = 9.835907521418337
    $ws.Range("L5").Value = 11.74350647444129
    $ws.Range("M5").Value = 21.95838308585659
    $ws.Range("N5").Value = 23.63692806849003
    $ws.Range("B6").Value = 28.0295845163205
    $ws.Range("C6").Value = 9.851803620306187
    $ws.Range("D6").Value = 3.709057880485368
    $ws.Range("F6").Value = 56.09132705366156
    $ws.Range("G6").Value = 3.800225665726589
    $ws.Range("J6").Value = 9.836240360504899
    $ws.Range("L6").Value = 11.74449593296603
    $ws.Range("M6").Value = 21.95786283727399
    $ws.Range("N6").Value = 23.63947408423449
    $ws.Range("B7").Value = 28.09602169251018
    $ws.Range("C7").Value = 9.962403961671312
    $ws.Range("D7").Value = 3.727919501686715
    $ws.Range("F7").Value = 56.1086440509082
    $ws.Range("G7").Value = 3.798845377647443
    $ws.Range("J7").Value = 9.833951444225736
    $ws.Range("L7").Value = 11.737760868314
    $ws.Range("M7").Value = 21.96188230424066
    $ws.Range("N7").Value = 23.62197874447249
    $ws.Range("B8").Value = 28.40490848097757
    $ws.Range("C8").Value = 10.44400769509147
    $ws.Range("D8").Value = 3.808740790689477
    $ws.Range("F8").Value = 56.21273266352
    $ws.Range("G8").Value = 3.793063676653002
    $ws.Range("J8").Value = 9.824378100595222
    $ws.Range("L8").Value = 11.71132636270364
    $ws.Range("M8").Value = 21.98970723447078
    $ws.Range("N8").Value = 23.54914948390467
    $ws.Range("B9").Value = 29.06313999488848
    $ws.Range("C9").Value = 11.35466862806599
    $ws.Range("D9").Value = 3.958178122253711
    $ws.Range("F9").Value = 56.51601020696141
    $ws.Range("G9").Value = 3.782827669338274
    $ws.Range("J9").Value = 9.807488023456786
    $ws.Range("L9").Value = 11.67128107738208
    $ws.Range("M9").Value = 22.08030894500214
    $ws.Range("N9").Value = 23.42204387815716
    $ws.Range("B10").Value = 29.57492025507288
    $ws.Range("C10").Value = 11.99502731106717
    $ws.Range("D10").Value = 4.06213322687553
    $ws.Range("F10").Value = 56.79881293847562
    $ws.Range("G10").Value = 3.775970070036327
    $ws.Range("J10").Value = 9.796216073330413
    $ws.Range("L10").Value = 11.64903324820778
    $ws.Range("M10").Value = 22.16858744706671
    $ws.Range("N10").Value = 23.33822962239908
    $ws.Range("B11").Value = 29.81299806449203
    $ws.Range("C11").Value = 12.27867702259136
    $ws.Range("D11").Value = 4.108096237586176
    $ws.Range("F11").Value = 56.94033876548721
    $ws.Range("G11").Value = 3.77299237611715
    $ws.Range("J11").Value = 9.791332789745223
    $ws.Range("L11").Value = 11.64046606023259
    $ws.Range("M11").Value = 22.21339177311807
    $ws.Range("N11").Value = 23.30217716679218
    $ws.Range("B12").Value = 29.90383460124243
    $ws.Range("C12").Value = 12.38488591900053
    $ws.Range("D12").Value = 4.125305597988416
    $ws.Range("F12").Value = 56.99576512413854
    $ws.Range("G12").Value = 3.771885052838776
    $ws.Range("J12").Value = 9.789518585972115
    $ws.Range("L12").Value = 11.63744491571997
    $ws.Range("M12").Value = 22.23101875033208
    $ws.Range("N12").Value = 23.28882347605137
    $ws.Range("B13").Value = 29.88424222468143
    $ws.Range("C13").Value = 12.36206694407065
    $ws.Range("D13").Value = 4.121608041401364
    $ws.Range("F13").Value = 56.98374683261521
    $ws.Range("G13").Value = 3.772122635580261
    $ws.Range("J13").Value = 9.789907753446531
    $ws.Range("L13").Value = 11.63808565784313
    $ws.Range("M13").Value = 22.22719321535806
    $ws.Range("N13").Value = 23.29168614878466
    $ws.Range("B14").Value = 29.82045801620747
    $ws.Range("C14").Value = 12.28743954820858
    $ws.Range("D14").Value = 4.10951601471941
    $ws.Range("F14").Value = 56.94486209046237
    $ws.Range("G14").Value = 3.772900870599395
    $ws.Range("J14").Value = 9.791182833729476
    $ws.Range("L14").Value = 11.64021304026787
    $ws.Range("M14").Value = 22.21482876261416
    $ws.Range("N14").Value = 23.30107256514391
    $ws.Range("B15").Value = 29.7814748566587
    $ws.Range("C15").Value = 12.24156860033137
    $ws.Range("D15").Value = 4.102083650918821
    $ws.Range("F15").Value = 56.92128226902895
    $ws.Range("G15").Value = 3.773380197096658
    $ws.Range("J15").Value = 9.791968409422731
    $ws.Range("L15").Value = 11.64154516288402
    $ws.Range("M15").Value = 22.20734096475024
    $ws.Range("N15").Value = 23.30686090810401
    $ws.Range("B16").Value = 29.55946083584973
    $ws.Range("C16").Value = 11.9763270837746
    $ws.Range("D16").Value = 4.059102300362508
    $ws.Range("F16").Value = 56.78982152634759
    $ws.Range("G16").Value = 3.776167512434114
    $ws.Range("J16").Value = 9.796540111962331
    $ws.Range("L16").Value = 11.64962436676468
    $ws.Range("M16").Value = 22.16575224162343
    $ws.Range("N16").Value = 23.34062750715593
    $ws.Range("B17").Value = 29.42455774023318
    $ws.Range("C17").Value = 11.81157293754594
    $ws.Range("D17").Value = 4.032391142327256
    $ws.Range("F17").Value = 56.71246074088248
    $ws.Range("G17").Value = 3.77791367839912
    $ws.Range("J17").Value = 9.799407186976641
    $ws.Range("L17").Value = 11.65497835446592
    $ws.Range("M17").Value = 22.14142374098401
    $ws.Range("N17").Value = 23.36187382853584
    $ws.Range("B18").Value = 29.34746534860467
    $ws.Range("C18").Value = 11.71609571556522
    $ws.Range("D18").Value = 4.016902999896208
    $ws.Range("F18").Value = 56.66917820298873
    $ws.Range("G18").Value = 3.77893138768043
    $ws.Range("J18").Value = 9.801079264203045
    $ws.Range("L18").Value = 11.6582040666711
    $ws.Range("M18").Value = 22.12786843464806
    $ws.Range("N18").Value = 23.37428945218527
    $ws.Range("B19").Value = 29.3214514070316
    $ws.Range("C19").Value = 11.68364930121321
    $ws.Range("D19").Value = 4.01163774152324
    $ws.Range("F19").Value = 56.65473238247638
    $ws.Range("G19").Value = 3.779278265399989
    $ws.Range("J19").Value = 9.801649357730769
    $ws.Range("L19").Value = 11.65932136364057
    $ws.Range("M19").Value = 22.12335425094766
    $ws.Range("N19").Value = 23.37852670665094
    $ws.Range("B20").Value = 29.43886717120623
    $ws.Range("C20").Value = 11.82918608006789
    $ws.Range("D20").Value = 4.035247522975948
    $ws.Range("F20").Value = 56.72057047591446
    $ws.Range("G20").Value = 3.777726414237591
    $ws.Range("J20").Value = 9.799099601451074
    $ws.Range("L20").Value = 11.65439328029556
    $ws.Range("M20").Value = 22.14396829253663
    $ws.Range("N20").Value = 23.35959190705693
    $ws.Range("B21").Value = 29.83917509184291
    $ws.Range("C21").Value = 12.30939283174444
    $ws.Range("D21").Value = 4.113073089758227
    $ws.Range("F21").Value = 56.9562338758808
    $ws.Range("G21").Value = 3.772671735202249
    $ws.Range("J21").Value = 9.790807363336812
    $ws.Range("L21").Value = 11.63958212529748
    $ws.Range("M21").Value = 22.21844263880021
    $ws.Range("N21").Value = 23.29830744167934
    $ws.Range("B22").Value = 30.10473701967654
    $ws.Range("C22").Value = 12.61618407832912
    $ws.Range("D22").Value = 4.162792994389029
    $ws.Range("F22").Value = 57.1209309354124
    $ws.Range("G22").Value = 3.769486266591958
    $ws.Range("J22").Value = 9.785591783709284
    $ws.Range("L22").Value = 11.63120219652038
    $ws.Range("M22").Value = 22.27096209463414
    $ws.Range("N22").Value = 23.25999503766701
    $ws.Range("B23").Value = 29.96266650946522
    $ws.Range("C23").Value = 12.45311979855966
    $ws.Range("D23").Value = 4.136362777404521
    $ws.Range("F23").Value = 57.03205868996785
    $ws.Range("G23").Value = 3.771175653436538
    $ws.Range("J23").Value = 9.788356831854408
    $ws.Range("L23").Value = 11.6355558838757
    $ws.Range("M23").Value = 22.24258223539817
    $ws.Range("N23").Value = 23.28028375511849
    $ws.Range("B24").Value = 29.43239642275629
    $ws.Range("C24").Value = 11.82122553100465
    $ws.Range("D24").Value = 4.033956562621077
    $ws.Range("F24").Value = 56.71690034678274
    $ws.Range("G24").Value = 3.777811033323624
    $ws.Range("J24").Value = 9.799238586856053
    $ws.Range("L24").Value = 11.65465733245914
    $ws.Range("M24").Value = 22.14281655627617
    $ws.Range("N24").Value = 23.36062293823021
    $ws.Range("B25").Value = 28.87983539805398
    $ws.Range("C25").Value = 11.11281931834974
    $ws.Range("D25").Value = 3.918754392649204
    $ws.Range("F25").Value = 56.42338732843267
    $ws.Range("G25").Value = 3.78547974194996
    $ws.Range("J25").Value = 9.811856754106879
    $ws.Range("L25").Value = 11.68085342662645
    $ws.Range("M25").Value = 22.05196398252703
    $ws.Range("N25").Value = 23.58581786079055
